# "added test for load flow calculation"
#
# The sheet had a trailing empty data row (row 10, just formula shells) right
# above the summary block (rows 11-13: sum[min]/sum[h]/sum[working weeks]).
# This change adds one more real data row of working hours, pushes the empty
# placeholder row down, and shifts the summary block down by one row,
# updating the running SUM() range to include the newly-added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 11 (and everything below, i.e. the old summary rows 11-13) down by
# one row. This leaves a fresh copy of the old "empty" row 10 template at the
# new row 11, and keeps row 10 free for the new data.
$ws.Rows.Item(11).Insert()

# Fill in the new working-hours entry on row 10.
$ws.Range("A10").Value = 2014
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 20
$ws.Range("D10").Value = 0.78125
$ws.Range("E10").Value = 0.90277777777777779
$ws.Range("F10").Formula = "=(E10-D10)*24*60"

# The "sum [min]" total (now on row 12) needs to cover the new data row plus
# the blank row 11 right above it, same pattern as before the insert.
$ws.Range("F12").Formula = "=SUM(F2:F11)"

# Selection moved to I10 in the saved workbook.
$ws.Range("I10").Select()
